# Modification of plan test
# - Fix wording in a couple of cells (typo / rewording corrections)
# - Update the current selection / view position on the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix "La page produite" -> "La page produit"
$ws.Range("B3").Value = "La page produit donne la description (de manière dynamique) du produit sélectionné en page d'accueil."

# Fix "description de produit" -> "description du produit"
$ws.Range("D3").Value = "Affiche la description du produit choisie."

# Reword the cart summary description and move it to D7
$ws.Range("D7").Value = " À la page panier le résumé des différentes informations des produits que l'utilisateur à choisis est affiché."

# Update the active selection / view to match the new state
$ws.Activate()
$ws.Range("C9").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 3
